# Dictionary.xlsx update: fit datamodel changes
#  - CollectionEvents sheet: drop the separate "ageMin" / "ageMax" columns
#    and replace them with a single "ageGroups" column (the old
#    "subcohorts" column shifts left to take column D).
#  - CollectionEvents becomes the active/selected sheet, with C4 selected.
#  - SourceTables is no longer the selected sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CollectionEvents")

# Remove the old "ageMax" column entirely (shifts "subcohorts" left into D)
# and turn the old "ageMin" header (now column C) into "ageGroups".
$ws.Range("D1").EntireColumn.Delete() | Out-Null
$ws.Range("C1").Value = "ageGroups"

# Match the narrower, best-fit-ish width the new header needs.
$ws.Columns.Item(3).ColumnWidth = 8.3

# Make CollectionEvents the active sheet with C4 selected, and drop the
# previous selection/active state on SourceTables.
$ws.Activate()
$ws.Range("C4").Select() | Out-Null
